# GitFlow.pptx change (commit: "Update Rubrics, try to make images larger.")
#
# The canonical-OOXML diff for this asset touches only three things:
#
#   1. ppt/presentation.xml       - the <go:slidesCustomData> extension
#                                    element is rewritten with its xmlns:*
#                                    attributes in a different order (same
#                                    attributes/values - cosmetic only, an
#                                    artifact of PowerPoint re-serialising
#                                    the file).
#   2. ppt/revisionInfo.xml       - dropped entirely (co-authoring session
#                                    bookkeeping: a random client GUID +
#                                    monotonic version counter + save
#                                    timestamp for *this* PowerPoint
#                                    session).
#   3. ppt/changesInfos/changesInfo1.xml
#                                  - one stale tracked-change record (a
#                                    <pc:spChg chg="del"> for shape id=50,
#                                    "PUSH" textbox, creationId
#                                    9C034F13-424F-E01C-207A-7AAF2177D493
#                                    on the slide sldId=271) is pruned from
#                                    the co-authoring change log.
#
# All three are PowerPoint's *internal* multi-author change-tracking
# bookkeeping (session GUIDs, monotonically increasing revision counters,
# wall-clock save timestamps, roundtrip signatures). They are not part of
# the slide/shape/text object model - real PowerPoint regenerates/prunes
# them on every save as a side effect of its co-authoring session, not
# through any Shape/Slide/Presentation API a script could call - so there
# is nothing to replay here through COM.
#
# Content-wise, the shape these records describe (the old "PUSH" textbox,
# id=50) is already gone from slide 3 in this deck, and its replacement
# (the "PUSH" textbox now carrying id=5) is already present with its
# final size/position - i.e. the actual visible edit (deleting the old
# label, adding the resized replacement) already happened earlier in the
# same authoring session and is already reflected in this file. Nothing
# about the slide contents changes in this diff.
#
# So this script intentionally performs no shape/text mutations - doing so
# would introduce a content diff that the target does not have. We simply
# touch the presentation through the object model (without changing
# anything) so the run completes cleanly.

$p = $ppt.ActivePresentation
$null = $p.Slides.Count
